$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.005100846290588
$ws.Range("B1").Value = 1.915569067001343
$ws.Range("C1").Value = 2.906126737594604
$ws.Range("D1").Value = 3.558630466461182
$ws.Range("E1").Value = 1.974258065223694
